# Adds a new "2022-Q3" quarter to the workbook:
#   1. Insert a new worksheet named "2022-Q3" right after "总计".
#   2. Insert a new row into the "总计" summary table for 2022-Q3
#      (pushing the existing quarters down by one row).
#   3. Populate the new "2022-Q3" worksheet with its fund-holding detail.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value to a cell as *text*, even if it looks numeric
# (e.g. fund codes like "010147"), without leaving a stray NumberFormat
# style behind on the cell. We stage the literal through a formula on a
# scratch cell (far outside any used range), copy its computed value over
# as a values-only paste, then clear the scratch cell.
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $row, $col, [string]$text) {
    $escaped = $text -replace '"', '""'
    $stage = $ws.Cells.Item(500, 500)
    $stage.Formula = '="' + $escaped + '"'
    [void]$stage.Copy()
    [void]$ws.Cells.Item($row, $col).PasteSpecial(-4163) # xlPasteValues
    [void]$stage.Clear()
}

# Helper: copy just the formatting (style) of a source cell onto a
# destination cell, then set its numeric value.
function Set-StyledNumber($destWs, $row, $col, $srcCell, $value) {
    [void]$srcCell.Copy()
    [void]$destWs.Cells.Item($row, $col).PasteSpecial(-4122) # xlPasteFormats
    $destWs.Cells.Item($row, $col).Value = $value
}

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($null, $totalSheet, 1)
$newSheet.Name = "2022-Q3"

# NOTE: re-fetch "2022-Q2" by name *after* inserting the new sheet — sheet
# variables here resolve by position, so grabbing it beforehand and holding
# onto the reference across the insert would silently alias the new sheet
# once everything shifts down a slot.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------------
# 2. Insert a new row 2 in "总计" for the 2022-Q3 summary figures.
# ---------------------------------------------------------------------------
[void]$totalSheet.Rows.Item(2).Insert()

Set-StyledNumber $totalSheet 2 1 $totalSheet.Cells.Item(3, 1) 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 9
$totalSheet.Cells.Item(2, 4).Value = 0.64

# The row-insert bleeds formatting onto B2:D2 (it copies row 1's style);
# clear that back to the unstyled look the other data rows use.
[void]$totalSheet.Range("B3:D3").Copy()
[void]$totalSheet.Range("B2:D2").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Populate the new "2022-Q3" worksheet.
# ---------------------------------------------------------------------------

# Header row (copy formatting from the equivalent header on "2022-Q2").
[void]$q2Sheet.Range("A1:H1").Copy()
[void]$newSheet.Range("A1:H1").PasteSpecial(-4163) # values
[void]$q2Sheet.Range("A1:H1").Copy()
[void]$newSheet.Range("A1:H1").PasteSpecial(-4122) # formats
$excel.CutCopyMode = 0

$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

$aStyleSrc = $q2Sheet.Cells.Item(2, 1)

$rows = @(
    @{ idx = 0; code = "010147"; name = "博道嘉兴一年持有期混合";     scale = "11.13"; pos = "90.11"; pct = "3.08"; mv = "0.3428"; rank = 10 },
    @{ idx = 1; code = "001305"; name = "九泰天富改革新动力混合A";   scale = "1.86";  pos = "94.71"; pct = "6.19"; mv = "0.1151"; rank = 10 },
    @{ idx = 2; code = "001844"; name = "九泰久益灵活配置混合C";     scale = "0.98";  pos = "93.32"; pct = "7.18"; mv = "0.0704"; rank = 9  },
    @{ idx = 3; code = "001782"; name = "九泰久益灵活配置混合A";     scale = "0.53";  pos = "93.32"; pct = "7.18"; mv = "0.0381"; rank = 9  },
    @{ idx = 4; code = "011686"; name = "创金合信先进装备股票C";     scale = "0.57";  pos = "80.17"; pct = "6.03"; mv = "0.0344"; rank = 7  },
    @{ idx = 5; code = "011685"; name = "创金合信先进装备股票A";     scale = "0.25";  pos = "80.17"; pct = "6.03"; mv = "0.0151"; rank = 7  },
    @{ idx = 6; code = "009912"; name = "九泰天富改革新动力混合C";   scale = "0.17";  pos = "94.71"; pct = "6.19"; mv = "0.0105"; rank = 10 },
    @{ idx = 7; code = "006182"; name = "格林伯锐灵活配置混合C";     scale = "0.14";  pos = "85.60"; pct = "6.19"; mv = "0.0087"; rank = 2  },
    @{ idx = 8; code = "006181"; name = "格林伯锐灵活配置混合A";     scale = "0.13";  pos = "85.60"; pct = "6.19"; mv = "0.0080"; rank = 2  }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $row = 2 + $i

    Set-StyledNumber $newSheet $row 1 $aStyleSrc $r.idx

    Set-TextValue $newSheet $row 2 $r.code
    Set-TextValue $newSheet $row 3 $r.name
    Set-TextValue $newSheet $row 4 $r.scale
    Set-TextValue $newSheet $row 5 $r.pos
    Set-TextValue $newSheet $row 6 $r.pct
    Set-TextValue $newSheet $row 7 $r.mv

    $newSheet.Cells.Item($row, 8).Value = $r.rank
}

$excel.CutCopyMode = 0
